# Automatic update of files.
#
# The 11 observation records in rows 2-12 of the "Artfynd" sheet get
# reshuffled: the same records stay, but each one moves to a different row.
# Concretely, new-row R ends up holding the data that used to live in
# old-row $mapping[R] (a permutation of 2..12 derived from the target diff).
#
# Only a subset of columns actually carries differing data across the moved
# rows (A, B, D, E, F, G, H, P, Q, R, Y, AA, AC, AW, AX); every other column
# in this row range is blank/absent for all 11 rows, so it is left alone -
# except for a few incidental "ghost" empty cells (J, N, AF) whose mere
# presence/absence shifts along with the row move; those are cleared where
# the destination row should end up without them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters -> 1-based column index, for the columns that hold real data.
$colIndex = @{
    "A"  = 1
    "B"  = 2
    "D"  = 4
    "E"  = 5
    "F"  = 6
    "G"  = 7
    "H"  = 8
    "P"  = 16
    "Q"  = 17
    "R"  = 18
    "Y"  = 25
    "AA" = 27
    "AC" = 29
    "AW" = 49
    "AX" = 50
}

$dataCols = @("A","B","D","E","F","G","H","P","Q","R","Y","AA","AC","AW","AX")

# Columns whose values look like dates ("2023-09-13" etc.) and must be kept
# as literal text instead of being auto-parsed into date serials by Excel.
$textProtectCols = @("Y","AA")

# new row number -> old row number that its data should come from.
$mapping = @{
    2  = 10
    3  = 8
    4  = 6
    5  = 7
    6  = 3
    7  = 9
    8  = 2
    9  = 11
    10 = 12
    11 = 4
    12 = 5
}

# 1) Snapshot every relevant cell from rows 2-12 before writing anything, so
#    that source rows being overwritten by other rows first can't clobber
#    data that is still needed later.
$snapshot = @{}
foreach ($r in 2..12) {
    $rowVals = @{}
    foreach ($col in $dataCols) {
        $idx = $colIndex[$col]
        $rowVals[$col] = $ws.Cells.Item($r, $idx).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshotted values back out in their permuted positions.
foreach ($newRow in 2..12) {
    $oldRow = $mapping[$newRow]
    $srcVals = $snapshot[$oldRow]
    foreach ($col in $dataCols) {
        $idx = $colIndex[$col]
        $val = $srcVals[$col]
        if (($textProtectCols -contains $col) -and ($val -ne $null)) {
            # Force text so Excel doesn't reinterpret e.g. "2023-09-13" as a date.
            $ws.Cells.Item($newRow, $idx).Value = "'" + $val
        } else {
            $ws.Cells.Item($newRow, $idx).Value = $val
        }
    }
}

# 3) A handful of incidental empty placeholder cells (J, N, AF) need to
#    disappear from specific destination rows to match the moved rows'
#    original shape.
$clearCells = @(
    @{ Row = 2;  Col = 10 },  # J2
    @{ Row = 2;  Col = 14 },  # N2
    @{ Row = 2;  Col = 32 },  # AF2
    @{ Row = 3;  Col = 10 },  # J3
    @{ Row = 3;  Col = 14 },  # N3
    @{ Row = 3;  Col = 32 },  # AF3
    @{ Row = 4;  Col = 14 },  # N4
    @{ Row = 5;  Col = 14 },  # N5
    @{ Row = 11; Col = 10 },  # J11
    @{ Row = 11; Col = 32 }   # AF11
)

foreach ($cc in $clearCells) {
    $ws.Cells.Item($cc.Row, $cc.Col).ClearContents()
}
